$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats")

# Row 2 / Row 8
$ws.Range("D2").Value = 0.0001182807609438896
$ws.Range("E2").Value = 0.02107867179438472
$ws.Range("G2").Value = 0.002035129815340042
$ws.Range("H2").Value = 0.003853811882436275
$ws.Range("I2").Value = 0.003970915451645851
$ws.Range("J2").Value = 0.007748536299914122
$ws.Range("K2").Value = 0.001373959239572287

$ws.Range("D8").Value = 0.0001182807609438896
$ws.Range("E8").Value = 0.02107867179438472
$ws.Range("G8").Value = 0.002035129815340042
$ws.Range("H8").Value = 0.003853811882436275
$ws.Range("I8").Value = 0.003970915451645851
$ws.Range("J8").Value = 0.007748536299914122
$ws.Range("K8").Value = 0.001373959239572287

# Row 3 / Row 9
$ws.Range("D3").Value = 0.0006056078709661961
$ws.Range("E3").Value = 0.01406853832304478
$ws.Range("G3").Value = 0.001336473971605301
$ws.Range("H3").Value = 0.003145434428006411
$ws.Range("I3").Value = 0.002672490198165178
$ws.Range("J3").Value = 0.005162716843187809
$ws.Range("K3").Value = 0.0005572894588112831

$ws.Range("D9").Value = 0.0006056078709661961
$ws.Range("E9").Value = 0.01406853832304478
$ws.Range("G9").Value = 0.001336473971605301
$ws.Range("H9").Value = 0.003145434428006411
$ws.Range("I9").Value = 0.002672490198165178
$ws.Range("J9").Value = 0.005162716843187809
$ws.Range("K9").Value = 0.0005572894588112831

# Row 4 / Row 10
$ws.Range("D4").Value = 0.0005990383215248585
$ws.Range("E4").Value = 0.01383312372490764
$ws.Range("G4").Value = 0.001268480438739061
$ws.Range("H4").Value = 0.003058651462197304
$ws.Range("I4").Value = 0.002586016431450844
$ws.Range("J4").Value = 0.004906681831926107
$ws.Range("K4").Value = 0.0006383815780282021

$ws.Range("D10").Value = 0.0005990383215248585
$ws.Range("E10").Value = 0.01383312372490764
$ws.Range("G10").Value = 0.001268480438739061
$ws.Range("H10").Value = 0.003058651462197304
$ws.Range("I10").Value = 0.002586016431450844
$ws.Range("J10").Value = 0.004906681831926107
$ws.Range("K10").Value = 0.0006383815780282021

# Row 5 / Row 11
$ws.Range("D5").Value = 0.0001673344522714615
$ws.Range("E5").Value = 0.02015670575201511
$ws.Range("G5").Value = 0.002040540799498558
$ws.Range("H5").Value = 0.003819882404059172
$ws.Range("I5").Value = 0.004338433500379324
$ws.Range("J5").Value = 0.007055392023175955
$ws.Range("K5").Value = 0.0008724918588995934

$ws.Range("D11").Value = 0.0001673344522714615
$ws.Range("E11").Value = 0.02015670575201511
$ws.Range("G11").Value = 0.002040540799498558
$ws.Range("H11").Value = 0.003819882404059172
$ws.Range("I11").Value = 0.004338433500379324
$ws.Range("J11").Value = 0.007055392023175955
$ws.Range("K11").Value = 0.0008724918588995934

# Row 6 / Row 12
$ws.Range("D6").Value = 0.001372099854052067
$ws.Range("E6").Value = 0.05951906181871891
$ws.Range("G6").Value = 0.003421592991799116
$ws.Range("H6").Value = 0.007727768737822771
$ws.Range("I6").Value = 0.03768706135451794
$ws.Range("J6").Value = 0.006112938746809959
$ws.Range("K6").Value = 0.001296191941946745

$ws.Range("D12").Value = 0.001372099854052067
$ws.Range("E12").Value = 0.05951906181871891
$ws.Range("G12").Value = 0.003421592991799116
$ws.Range("H12").Value = 0.007727768737822771
$ws.Range("I12").Value = 0.03768706135451794
$ws.Range("J12").Value = 0.006112938746809959
$ws.Range("K12").Value = 0.001296191941946745
